$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $saveStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $saveStyle
}

Set-TextValue $ws "D2" "299.16"
Set-TextValue $ws "E2" "-1.30%"
Set-TextValue $ws "D3" "31.53"
Set-TextValue $ws "E3" "-0.57%"
Set-TextValue $ws "D4" "5.095"
Set-TextValue $ws "E4" "-1.22%"
Set-TextValue $ws "D5" "0.07865"
Set-TextValue $ws "E5" "0.67%"
Set-TextValue $ws "D6" "2.277"
Set-TextValue $ws "E6" "-1.41%"
Set-TextValue $ws "D7" "7.806"
Set-TextValue $ws "E7" "-1.80%"
Set-TextValue $ws "D8" "3.853"
Set-TextValue $ws "E8" "-0.47%"
Set-TextValue $ws "D9" "0.9206"
Set-TextValue $ws "E9" "1.39%"
Set-TextValue $ws "D10" "0.1749"
Set-TextValue $ws "E10" "1.13%"
Set-TextValue $ws "D11" "0.07567"
Set-TextValue $ws "E11" "3.05%"
Set-TextValue $ws "D12" "0.09298"
Set-TextValue $ws "E12" "14.44%"
Set-TextValue $ws "D13" "0.03009"
Set-TextValue $ws "E13" "-0.27%"
Set-TextValue $ws "D14" "0.1002"
Set-TextValue $ws "E14" "0.79%"
Set-TextValue $ws "D15" "0.001509"
Set-TextValue $ws "E15" "-0.38%"
Set-TextValue $ws "D16" "0.006100"
Set-TextValue $ws "E16" "1.43%"
Set-TextValue $ws "D17" "3.473"
Set-TextValue $ws "E17" "-0.70%"
Set-TextValue $ws "E18" "0.15%"
Set-TextValue $ws "E19" "0.86%"
Set-TextValue $ws "D20" "0.1309"
Set-TextValue $ws "E20" "-2.07%"
Set-TextValue $ws "D21" "3.935"
Set-TextValue $ws "E21" "-15.73%"
Set-TextValue $ws "D22" "0.1712"
Set-TextValue $ws "E22" "9.44%"
Set-TextValue $ws "D23" "0.04621"
Set-TextValue $ws "E23" "-0.72%"
Set-TextValue $ws "D24" "0.001255"
Set-TextValue $ws "E24" "-0.52%"
Set-TextValue $ws "D25" "0.004470"
Set-TextValue $ws "E25" "-1.10%"
Set-TextValue $ws "D26" "0.0001251"
Set-TextValue $ws "E26" "-7.28%"
Set-TextValue $ws "D27" "0.0003400"
Set-TextValue $ws "E27" "24.11%"
Set-TextValue $ws "D39" "0.01738"
Set-TextValue $ws "E39" "-3.41%"
Set-TextValue $ws "D40" "0.04601"
Set-TextValue $ws "E40" "0.84%"
Set-TextValue $ws "D41" "0.006927"
Set-TextValue $ws "E41" "-4.97%"
Set-TextValue $ws "D42" "0.1359"
Set-TextValue $ws "D43" "0.002192"
Set-TextValue $ws "E43" "-2.10%"
Set-TextValue $ws "D44" "0.009762"
Set-TextValue $ws "E44" "-9.15%"
Set-TextValue $ws "D45" "0.00006290"
Set-TextValue $ws "E45" "-3.03%"
Set-TextValue $ws "E46" "0.18%"
Set-TextValue $ws "D47" "0.007985"
Set-TextValue $ws "E47" "-19.32%"
Set-TextValue $ws "D48" "1.154"
Set-TextValue $ws "E48" "40.68%"
Set-TextValue $ws "D49" "0.00002101"
Set-TextValue $ws "E49" "0.18%"
Set-TextValue $ws "D50" "0.0002001"
Set-TextValue $ws "E50" "0.18%"
